# Auto-generated script to apply scheduled market-data refresh to Sheets workbook
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 628.5714
$ws.Range("I18").Value = 240
$ws.Range("J18").Value = 1600
$ws.Range("K18").Value = 240
$ws.Range("L18").Value = 1600
$ws.Range("M18").Value = 44
$ws.Range("N18").Value = -2168
$ws.Range("H74").Value = 9016.75
$ws.Range("I74").Value = 3525.2727
$ws.Range("K74").Value = 3525.2727
$ws.Range("M74").Value = -2589.2727
$ws.Range("H76").Value = 66740784
$ws.Range("I76").Value = 79413.5
$ws.Range("K76").Value = 79413.5
$ws.Range("M76").Value = -79098.5
$ws.Range("H77").Value = 9016.75
$ws.Range("I77").Value = 3525.2727
$ws.Range("K77").Value = 17626.3635
$ws.Range("M77").Value = -12946.3635
$ws.Range("H79").Value = 66740784
$ws.Range("I79").Value = 79413.5
$ws.Range("K79").Value = 79413.5
$ws.Range("M79").Value = -78321.5
$ws.Range("H86").Value = 4388478
$ws.Range("J86").Value = 8773638
$ws.Range("L86").Value = 8773638
$ws.Range("N86").Value = -8775884
$ws.Range("H89").Value = 4388478
$ws.Range("J89").Value = 8773638
$ws.Range("L89").Value = 43868190
$ws.Range("N89").Value = -43879422
$ws.Range("H112").Value = 4343.3213
$ws.Range("J112").Value = 4716.52
$ws.Range("L112").Value = 14149.56
$ws.Range("N112").Value = -16365.56
$ws.Range("H116").Value = 4511.4443
$ws.Range("I116").Value = 4184.0557
$ws.Range("J116").Value = 5166.222
$ws.Range("K116").Value = 4184.0557
$ws.Range("L116").Value = 5166.222
$ws.Range("M116").Value = -742.0556999999999
$ws.Range("N116").Value = -12050.222
$ws.Range("H132").Value = 2668.0625
$ws.Range("I132").Value = 2049.2144
$ws.Range("K132").Value = 6147.6432
$ws.Range("M132").Value = -3617.6432
$ws.Range("H138").Value = 6552.406
$ws.Range("I138").Value = 5882.3
$ws.Range("J138").Value = 6665.983
$ws.Range("K138").Value = 17646.9
$ws.Range("L138").Value = 19997.949
$ws.Range("M138").Value = -12506.9
$ws.Range("N138").Value = -30277.949

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2996.4265
$ws.Range("I32").Value = 2110.8196
$ws.Range("K32").Value = 2110.8196
$ws.Range("M32").Value = -1823.8196
$ws.Range("H61").Value = 3807.9312
$ws.Range("I61").Value = 2364.2104
$ws.Range("K61").Value = 2364.2104
$ws.Range("M61").Value = -2152.2104
$ws.Range("H76").Value = 47000
$ws.Range("J76").Value = 47000
$ws.Range("L76").Value = 47000
$ws.Range("N76").Value = -47676
$ws.Range("H79").Value = 47000
$ws.Range("J79").Value = 47000
$ws.Range("L79").Value = 47000
$ws.Range("N79").Value = -49340
$ws.Range("H110").Value = 417890.4
$ws.Range("J110").Value = 2013
$ws.Range("L110").Value = 2013
$ws.Range("N110").Value = -6103
$ws.Range("H122").Value = 5507
$ws.Range("I122").Value = 2999.5
$ws.Range("K122").Value = 8998.5
$ws.Range("M122").Value = -6548.5
$ws.Range("H132").Value = 3607.7673
$ws.Range("I132").Value = 1711.5518
$ws.Range("J132").Value = 7535.643
$ws.Range("K132").Value = 5134.6554
$ws.Range("L132").Value = 22606.929
$ws.Range("M132").Value = -2604.6554
$ws.Range("N132").Value = -27666.929
$ws.Range("H136").Value = 3807.9312
$ws.Range("I136").Value = 2364.2104
$ws.Range("K136").Value = 7092.6312
$ws.Range("M136").Value = -4542.6312

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1993.125
$ws.Range("I20").Value = 1354
$ws.Range("J20").Value = 2814.8572
$ws.Range("K20").Value = 1354
$ws.Range("L20").Value = 2814.8572
$ws.Range("M20").Value = -1107
$ws.Range("N20").Value = -3308.8572
$ws.Range("H25").Value = 1576.6666
$ws.Range("I25").Value = 865
$ws.Range("J25").Value = 3000
$ws.Range("K25").Value = 865
$ws.Range("L25").Value = 3000
$ws.Range("M25").Value = -630
$ws.Range("N25").Value = -3470
$ws.Range("H94").Value = 1740.9706
$ws.Range("I94").Value = 599.5
$ws.Range("K94").Value = 599.5
$ws.Range("M94").Value = -148.5
$ws.Range("H99").Value = 3762.625
$ws.Range("I99").Value = 3457.7778
$ws.Range("J99").Value = 4154.5713
$ws.Range("K99").Value = 3457.7778
$ws.Range("L99").Value = 4154.5713
$ws.Range("M99").Value = -1959.7778
$ws.Range("N99").Value = -7150.5713
$ws.Range("H105").Value = 39448.77
$ws.Range("I105").Value = 39448.77
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 39448.77
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -37701.77
$ws.Range("N105").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 361592.47
$ws.Range("I31").Value = 590577.3
$ws.Range("J31").Value = 7706.8184
$ws.Range("K31").Value = 590577.3
$ws.Range("L31").Value = 7706.8184
$ws.Range("M31").Value = -590282.3
$ws.Range("N31").Value = -8296.8184
$ws.Range("H34").Value = 361592.47
$ws.Range("I34").Value = 590577.3
$ws.Range("J34").Value = 7706.8184
$ws.Range("K34").Value = 590577.3
$ws.Range("L34").Value = 7706.8184
$ws.Range("M34").Value = -590375.3
$ws.Range("N34").Value = -8110.8184
$ws.Range("H58").Value = 211958
$ws.Range("I58").Value = 324877.22
$ws.Range("K58").Value = 324877.22
$ws.Range("M58").Value = -324674.22
$ws.Range("H132").Value = 3280.3
$ws.Range("I132").Value = 1789.174
$ws.Range("J132").Value = 5297.706
$ws.Range("K132").Value = 5367.522
$ws.Range("L132").Value = 15893.118
$ws.Range("M132").Value = -2837.522
$ws.Range("N132").Value = -20953.118
$ws.Range("H134").Value = 236654.47
$ws.Range("I134").Value = 2746.1853
$ws.Range("K134").Value = 8238.555899999999
$ws.Range("M134").Value = -5703.555899999999
$ws.Range("H136").Value = 211958
$ws.Range("I136").Value = 324877.22
$ws.Range("K136").Value = 974631.6599999999
$ws.Range("M136").Value = -972081.6599999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 93907.41
$ws.Range("I107").Value = 1068.5555
$ws.Range("J107").Value = 158180.47
$ws.Range("K107").Value = 3205.6665
$ws.Range("L107").Value = 474541.41
$ws.Range("M107").Value = -1285.6665
$ws.Range("N107").Value = -478381.41

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H60").Value = 5000
$ws.Range("I60").Value = 5000
$ws.Range("K60").Value = 5000
$ws.Range("M60").Value = -4425
$ws.Range("H80").Value = 955778.7
$ws.Range("I80").Value = 559257.8
$ws.Range("J80").Value = 3334903.8
$ws.Range("K80").Value = 559257.8
$ws.Range("L80").Value = 3334903.8
$ws.Range("M80").Value = -558259.8
$ws.Range("N80").Value = -3336899.8
$ws.Range("H83").Value = 955778.7
$ws.Range("I83").Value = 559257.8
$ws.Range("J83").Value = 3334903.8
$ws.Range("K83").Value = 2796289
$ws.Range("L83").Value = 16674519
$ws.Range("M83").Value = -2791297
$ws.Range("N83").Value = -16684503
$ws.Range("H113").Value = 2003582
$ws.Range("I113").Value = 5001455.5
$ws.Range("J113").Value = 4999.6665
$ws.Range("K113").Value = 5001455.5
$ws.Range("L113").Value = 4999.6665
$ws.Range("M113").Value = -4999285.5
$ws.Range("N113").Value = -9339.666499999999
$ws.Range("H132").Value = 426269.3
$ws.Range("I132").Value = 479267.25
$ws.Range("K132").Value = 1437801.75
$ws.Range("M132").Value = -1435271.75

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H33").Value = 14500
$ws.Range("I33").Value = 14500
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 14500
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -14210
$ws.Range("N33").ClearContents()
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("M50").ClearContents()
$ws.Range("H82").Value = 6300.1665
$ws.Range("I82").Value = 5560.2
$ws.Range("K82").Value = 5560.2
$ws.Range("M82").Value = -5199.2
$ws.Range("H85").Value = 6300.1665
$ws.Range("I85").Value = 5560.2
$ws.Range("K85").Value = 5560.2
$ws.Range("M85").Value = -4312.2
$ws.Range("H98").Value = 19955
$ws.Range("J98").Value = 19955
$ws.Range("L98").Value = 19955
$ws.Range("N98").Value = -25945
$ws.Range("H100").Value = 1644.4445
$ws.Range("I100").Value = 1435.8334
$ws.Range("J100").Value = 2061.6667
$ws.Range("K100").Value = 1435.8334
$ws.Range("L100").Value = 2061.6667
$ws.Range("M100").Value = -894.8334
$ws.Range("N100").Value = -3143.6667
$ws.Range("H136").Value = 4258.2964
$ws.Range("I136").Value = 3089.5789
$ws.Range("J136").Value = 7034
$ws.Range("K136").Value = 9268.736699999999
$ws.Range("L136").Value = 21102
$ws.Range("M136").Value = -6718.736699999999
$ws.Range("N136").Value = -26202

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("H49").Value = 3333333
$ws.Range("I49").Value = 3333333
$ws.Range("K49").Value = 3333333
$ws.Range("M49").Value = -3333103
$ws.Range("H132").Value = 51136
$ws.Range("J132").Value = 117460.89
$ws.Range("L132").Value = 352382.67
$ws.Range("N132").Value = -357442.67
$ws.Range("H136").Value = 67057.73
$ws.Range("I136").Value = 14223.031
$ws.Range("J136").Value = 405199.8
$ws.Range("K136").Value = 42669.093
$ws.Range("L136").Value = 1215599.4
$ws.Range("M136").Value = -40119.093
$ws.Range("N136").Value = -1220699.4
